$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" (column D) and "Volume(1h)" (column E) cells to the
# latest scraped values. Column D values are forced to remain plain
# text (matching the original inlineStr cells) via a temporary "@"
# number format so Excel does not reinterpret numeric-looking strings
# (e.g. "313.92", "0.4330") as actual numbers, which would silently
# drop significant trailing zeros / precision. The format is reset
# back to "Normal" style immediately after so no styling changes leak
# into the saved workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.736.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4330"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008986"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.605.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.980"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "120.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.248"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08903"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01930"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.844"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5084"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1656"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.318"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
